$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the TurretAlert status as Completed (turret alert sound implemented)
$ws.Range("E6").Value = "Completed"

# Mark the Jump status as Completed
$ws.Range("E17").Value = "Completed"

# Update the Ambience assets required text: drop the "(x3)" qualifier on lava bubbling
$ws.Range("D25").Value = "Wind loop SFX (x2), Lava bubbling SFX"

# Move the active selection to E7
$ws.Range("E7").Select()
